$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update the Tuesday ("SELASA") period 3-4 block (row 4-5) in column E:
# swap "B. ING" for "MATE" so it lines up with the new pairing and keep the
# merged (no-bottom-border) look shared by the row-4/row-5 pair.
$ws.Range("E4").Value = "MATE"
$ws.Range("E5").Value = "MATE"
$ws.Range("E5").Copy() | Out-Null
$ws.Range("E4").PasteSpecial(-4122) | Out-Null

# Update the Thursday ("KAMIS") period 7-8 block (row 9-10) in column E:
# swap "MATE" for "B. ING" and give row 10 the same full-border look as row 9.
$ws.Range("E9").Value = "B. ING"
$ws.Range("E10").Value = "B. ING"
$ws.Range("E9").Copy() | Out-Null
$ws.Range("E10").PasteSpecial(-4122) | Out-Null

$excel.CutCopyMode = 0

# Move the active selection to E5, matching the saved cursor position.
$ws.Range("E5").Select() | Out-Null
